$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.628.69"
$ws.Range("E2").Value = "  -1.77%  "

$ws.Range("D3").Value = "1.588.52"
$ws.Range("E3").Value = "  -2.32%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.98"
$ws.Range("E5").Value = "  -1.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.510"
$ws.Range("E6").Value = "  -2.66%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  -2.43%  "

$ws.Range("E9").Value = "  -2.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.57"
$ws.Range("E10").Value = "  -3.84%  "

$ws.Range("E11").Value = "  -1.75%  "

$ws.Range("D12").Value = "1.811.68"
$ws.Range("E12").Value = "  -2.28%  "

$ws.Range("D13").Value = "1.627.86"
$ws.Range("E13").Value = "  -0.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.02"
$ws.Range("E14").Value = "  -2.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.521"
$ws.Range("E15").Value = "  -4.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.74"
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("D17").Value = "26.613.50"
$ws.Range("E17").Value = "  -1.86%  "

$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  -2.28%  "

$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.00"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "207.97"
$ws.Range("E20").Value = "  -4.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.73"
$ws.Range("E21").Value = "  -3.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.24"
$ws.Range("E22").Value = "  -2.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.33"
$ws.Range("E23").Value = "  -4.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.86"
$ws.Range("E24").Value = "  -2.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.67"
$ws.Range("E25").Value = "  -0.97%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.24"
$ws.Range("E27").Value = "  -0.87%  "

$ws.Range("E28").Value = "  -3.46%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.28"
$ws.Range("E29").Value = "  -2.21%  "

$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("E31").Value = "  -1.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  -4.21%  "

$ws.Range("E33").Value = "  +20.31%  "

$ws.Range("E34").Value = "  -2.93%  "

$ws.Range("D35").Value = "1.306.59"
$ws.Range("E35").Value = "  -2.83%  "

$ws.Range("E36").Value = "  -1.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.48"
$ws.Range("E37").Value = "  -5.58%  "

$ws.Range("E38").Value = "  -3.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.828"
$ws.Range("E39").Value = "  -3.32%  "

$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.792"
$ws.Range("E41").Value = "  -1.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.36"
$ws.Range("E42").Value = "  +2.44%  "

$ws.Range("E43").Value = "  -3.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.58"
$ws.Range("E44").Value = "  -4.48%  "

$ws.Range("D45").Value = "1.724.88"
$ws.Range("E45").Value = "  -2.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.54"

$ws.Range("E47").Value = "  -0.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.836"
$ws.Range("E48").Value = "  -2.24%  "

$ws.Range("E49").Value = "  -1.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0977"
$ws.Range("E50").Value = "  -1.73%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.48"
$ws.Range("E51").Value = "  -1.68%  "
